$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cid values in column A (rows 11-18) as part of the Tube Data split cleanup
$ws.Range("A11").Value = 6100
$ws.Range("A12").Value = 6101
$ws.Range("A13").Value = 6102
$ws.Range("A14").Value = 6103
$ws.Range("A15").Value = 6104
$ws.Range("A16").Value = 6105
$ws.Range("A17").Value = 6106
$ws.Range("A18").Value = 6107

# Scroll the view and update the active selection
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
